# example1/settings.xlsx — "Add files via upload" commit.
#
# Net content-level change (row/column values were otherwise untouched):
#   1) A handful of descriptive strings in column B had spelling fixes:
#        "firts" -> "first", "rapresented" -> "represented", "set ot TRUE" -> "set to TRUE"
#   2) A new setting row (CPU / number of threads...) was appended as row 31.
#   3) The active selection moved to B31 (the newly added cell).
#
# (Cosmetic-only bytes that a real Excel re-save also perturbs -- fileVersion/
# revisionPtr guids, theme display name, exact column-width floats, the
# x14ac dyDescent default-row-height hint -- are artifacts of which Excel
# build wrote the file and are not reproducible/meaningful host actions, so
# they are intentionally left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value  = "single or paired if alignment or trimming is set to TRUE"
$ws.Range("B3").Value  = "adapter sequence to be trimmed  (if in paired mode applied to the first strand) if trimming is set to TRUE"
$ws.Range("B5").Value  = "number of bases to be clipped at 5' end (if in paired mode applied to the first strand) if trimming is set to TRUE"
$ws.Range("B7").Value  = "number of bases to be clipped at 3' end (if in paired mode applied to the first strand) if trimming is set to TRUE"
$ws.Range("B17").Value = "name of genes that must be represented in boxplots (1 in each column)"
$ws.Range("B18").Value = "name of genes that must be represented in heatmaps (1 in each column)"
$ws.Range("B20").Value = "number of top gene ontology terms represented"
$ws.Range("B23").Value = "number of PCA components to be represented (it must be less or equal to the number of samples)"
$ws.Range("B24").Value = "gene ontology terms (of the corresponding GO_type in the form GO:number) that must be represented in the enrichment plots (1 in each column)"
$ws.Range("B25").Value = "pathways names (with the official pathway id) that must be represented in pathview (1 in each column)"

$ws.Range("A31").Value = "CPU"
$ws.Range("B31").Value = "number of threads to be used for STAR indexing and alignment"
$ws.Range("C31").Value = 12

$ws.Range("B31").Select()
